# Add a new "Math" progress entry as row 15 of Table1 (A1:E14 -> A1:E15),
# mirroring the formatting of the row above it, then resize the table and
# leave the selection where Excel would land after typing the row (B16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from row 13's formatting (border/font already match the table body
# style used for the new row) and paste it onto row 15 before typing values.
$ws.Range("A13:E13").Copy()
$ws.Range("A15").PasteSpecial(-4122)

# Enter the new row's data.
$ws.Range("A15").Value = "Math"
$ws.Range("B15").Value = 45110
$ws.Range("C15").Value = 849.36
$ws.Range("D15").Value = 2524
$ws.Range("E15").Formula = "=IF(ROW()>2,(`$D`$2-D15)/`$D`$2,""NA"")"

# Match the exact style variants Excel produced for this row (the new
# border/fill/number-format combinations it derives on data entry).
$ws.Range("C15:D15").Interior.ColorIndex = -4142
$ws.Range("B15").Interior.ColorIndex = -4142
$ws.Range("E15").Interior.ColorIndex = -4142

# Grow Table1 so the new row becomes part of the table.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E15"))

# Final selection, as left behind after typing into B15 and pressing Enter.
$ws.Range("B16").Select()
